# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns for the
# 84f8000b-cf6b-4b2a-8656-2428c2200111 row (row 6) on both the zh-cn and
# de-de sheets, widens the Error Detail column, and links the new
# "Latest Target File" cell back to the source markdown on GitHub.

$wb = $excel.ActiveWorkbook

# Cornflower-blue (FF6495ED) used by the other hyperlink cells on these
# sheets, expressed as a VBA-style BGR-packed long for Font.Color.
$hyperlinkColor = 15570276

function Set-HandbackReportRow6 {
    param(
        $ws,
        [string]$targetFileDisplay,
        [string]$targetFileUrl,
        [string]$handbackFile,
        [string]$handbackDateTime,
        [string]$errorDetail
    )

    # I6 - Latest Target File: becomes a hyperlink to the source markdown.
    $ws.Hyperlinks.Add($ws.Range("I6"), $targetFileUrl, "", "", $targetFileDisplay)
    $ws.Range("I6").Font.Underline = $true
    $ws.Range("I6").Font.Color = $hyperlinkColor

    # J6 - Latest Handback File
    $ws.Range("J6").Value = $handbackFile

    # K6 - Latest Handback DateTime
    $ws.Range("K6").Value = $handbackDateTime

    # P6 - Error Detail
    $ws.Range("P6").Value = $errorDetail

    # Widen column P (Error Detail) so the long message is readable.
    $ws.Columns.Item(16).ColumnWidth = 39.166666666666664
}

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/901c950ad49bbe9cdc31081a039fd32840aef048/e2e/84f8000b-cf6b-4b2a-8656-2428c2200111.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b7c200ac9842fbee428563245f85417fcdf853e/e2e/84f8000b-cf6b-4b2a-8656-2428c2200111.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/901c950ad49bbe9cdc31081a039fd32840aef048/e2e/84f8000b-cf6b-4b2a-8656-2428c2200111.md."

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackReportRow6 -ws $wsZhCn `
    -targetFileDisplay "84f8000b-cf6b-4b2a-8656-2428c2200111.md" `
    -targetFileUrl $latestMdUrl `
    -handbackFile "84f8000b-cf6b-4b2a-8656-2428c2200111.c4b6d7089d29c93a17cb924f68327778ad2fc54d.zh-cn.xlf" `
    -handbackDateTime "2016-08-26 16:43:51" `
    -errorDetail $errorDetail

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackReportRow6 -ws $wsDeDe `
    -targetFileDisplay "84f8000b-cf6b-4b2a-8656-2428c2200111.md" `
    -targetFileUrl $latestMdUrl `
    -handbackFile "84f8000b-cf6b-4b2a-8656-2428c2200111.c4b6d7089d29c93a17cb924f68327778ad2fc54d.de-de.xlf" `
    -handbackDateTime "2016-08-26 16:43:58" `
    -errorDetail $errorDetail

Write-Output "Handback report updated for zh-cn and de-de (row 6)."
